$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '61.853.31'
Set-TextValue $ws.Range("E2") '  +1.20%  '
Set-TextValue $ws.Range("D3") '3.414.66'
Set-TextValue $ws.Range("E3") '  +3.77%  '
Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.00%  '
Set-TextValue $ws.Range("D5") '577.33'
Set-TextValue $ws.Range("E5") '  +2.36%  '
Set-TextValue $ws.Range("D6") '139.43'
Set-TextValue $ws.Range("E6") '  +9.75%  '
Set-TextValue $ws.Range("D7") '1.00'
Set-TextValue $ws.Range("E7") '  -0.11%  '
Set-TextValue $ws.Range("D8") '3.414.69'
Set-TextValue $ws.Range("E8") '  +3.89%  '
Set-TextValue $ws.Range("E9") '  +0.38%  '
Set-TextValue $ws.Range("D10") '7.69'
Set-TextValue $ws.Range("D11") '0.127'
Set-TextValue $ws.Range("E11") '  +8.53%  '
Set-TextValue $ws.Range("E12") '  +6.09%  '
Set-TextValue $ws.Range("D13") '3.998.72'
Set-TextValue $ws.Range("E13") '  +3.34%  '
Set-TextValue $ws.Range("E14") '  +2.06%  '
Set-TextValue $ws.Range("D15") '0.0000181'
Set-TextValue $ws.Range("E15") '  +8.82%  '
Set-TextValue $ws.Range("E16") '  +3.41%  '
Set-TextValue $ws.Range("D17") '25.71'
Set-TextValue $ws.Range("E17") '  +6.89%  '
Set-TextValue $ws.Range("D18") '61.926.03'
Set-TextValue $ws.Range("E18") '  +1.13%  '
Set-TextValue $ws.Range("B19") 'Polkadot'
Set-TextValue $ws.Range("C19") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D19") '5.96'
Set-TextValue $ws.Range("E19") '  +6.79%  '
Set-TextValue $ws.Range("B20") 'Chainlink'
Set-TextValue $ws.Range("C20") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D20") '14.13'
Set-TextValue $ws.Range("E20") '  +7.12%  '
Set-TextValue $ws.Range("E21") '  +6.65%  '
Set-TextValue $ws.Range("D22") '392.44'
Set-TextValue $ws.Range("E22") '  +11.48%  '
Set-TextValue $ws.Range("E23") '  +4.60%  '
Set-TextValue $ws.Range("D24") '3.556.72'
Set-TextValue $ws.Range("E24") '  +3.68%  '
Set-TextValue $ws.Range("B25") 'Dai'
Set-TextValue $ws.Range("C25") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D25") '1.00'
Set-TextValue $ws.Range("E25") '  +0.19%  '
Set-TextValue $ws.Range("B26") 'PEPE'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D26") '0.0000127'
Set-TextValue $ws.Range("E26") '  +19.31%  '
Set-TextValue $ws.Range("D27") '71.34'
Set-TextValue $ws.Range("E27") '  +3.43%  '
Set-TextValue $ws.Range("E28") '  +15.99%  '
Set-TextValue $ws.Range("D29") '7.82'
Set-TextValue $ws.Range("E29") '  +10.41%  '
Set-TextValue $ws.Range("D30") '0.999'
Set-TextValue $ws.Range("E30") '  -0.05%  '
Set-TextValue $ws.Range("D31") '8.34'
Set-TextValue $ws.Range("E31") '  +7.81%  '
Set-TextValue $ws.Range("E32") '  +7.68%  '
Set-TextValue $ws.Range("E33") '  +3.06%  '
Set-TextValue $ws.Range("D34") '3.449.88'
Set-TextValue $ws.Range("E34") '  +3.72%  '
Set-TextValue $ws.Range("E35") '  -0.05%  '
Set-TextValue $ws.Range("D36") '23.69'
Set-TextValue $ws.Range("E36") '  +5.27%  '
Set-TextValue $ws.Range("D37") '5.53'
Set-TextValue $ws.Range("E37") '  +6.01%  '
Set-TextValue $ws.Range("D38") '7.08'
Set-TextValue $ws.Range("E38") '  +5.08%  '
Set-TextValue $ws.Range("E39") '  +6.77%  '
Set-TextValue $ws.Range("D40") '162.00'
Set-TextValue $ws.Range("E40") '  -0.69%  '
Set-TextValue $ws.Range("D41") '0.0803'
Set-TextValue $ws.Range("E41") '  +7.26%  '
Set-TextValue $ws.Range("E42") '  +12.51%  '
Set-TextValue $ws.Range("E43") '  -0.01%  '
Set-TextValue $ws.Range("E44") '  +3.23%  '
Set-TextValue $ws.Range("B45") 'ONDO'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range("D45") '1.23'
Set-TextValue $ws.Range("E45") '  +10.04%  '
Set-TextValue $ws.Range("B46") 'Mantle'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D46") '0.776'
Set-TextValue $ws.Range("E46") '  +4.97%  '
Set-TextValue $ws.Range("D47") '41.26'
Set-TextValue $ws.Range("E47") '  +0.46%  '
Set-TextValue $ws.Range("D48") '23.44'
Set-TextValue $ws.Range("E48") '  +5.84%  '
Set-TextValue $ws.Range("D49") '7.03'
Set-TextValue $ws.Range("E49") '  +5.76%  '
Set-TextValue $ws.Range("E50") '  +9.34%  '
Set-TextValue $ws.Range("D51") '2.362.82'
Set-TextValue $ws.Range("E51") '  +10.24%  '
